$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.885.75'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.517.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.81%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.602'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.518.37'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("E10").Value = '  -1.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.94'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.425'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.132.24'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '30.54'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.57%  '
$ws.Range("E15").Value = '  -2.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.917.42'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("E17").Value = '  -2.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.521.21'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.90%  '
$ws.Range("E19").Value = '  -3.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '383.45'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.45%  '
$ws.Range("E22").Value = '  -0.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.550'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.73%  '
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.53'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.10%  '
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("E27").Value = '  -0.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.89'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.30%  '
$ws.Range("E29").Value = '  -1.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '24.71'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.91'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.22%  '
$ws.Range("E34").Value = '  -5.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.29'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.58'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '30.02'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +13.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '161.27'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.896'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.12%  '
$ws.Range("E41").Value = '  -4.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.61'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.29%  '
$ws.Range("E44").Value = '  -7.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.740.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0706'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.28'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.68%  '
$ws.Range("E48").Value = '  -2.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0298'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '324.71'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.02'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.19%  '
